$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Fill in the two new cells on row 3 (Bollinger band fix produced trailing PriceChange/UpDown values)
$ws.Range("X3").Value = -0.29999499999999557
$ws.Range("Y3").Value = "Down"

# Duplicate row 3 formatting down into row 4 by copy/paste so styles (date format, percent format) carry over
$ws.Range("A3:W3").Copy()
$ws.Range("A4").PasteSpecial(-4122)
$excel.CutCopyMode = 0

# Now populate row 4 with its own data values
$ws.Range("A4").Value = 42635.817361111112
$ws.Range("B4").Value = 11
$ws.Range("C4").Value = "Buy"
$ws.Range("D4").Value = 64
$ws.Range("E4").Value = 7650
$ws.Range("F4").Value = 454
$ws.Range("G4").Value = 75
$ws.Range("H4").Value = 23
$ws.Range("I4").Value = 100
$ws.Range("J4").Value = 0
$ws.Range("K4").Value = 12332
$ws.Range("L4").Value = 82
$ws.Range("M4").Value = 26
$ws.Range("N4").Value = 16
$ws.Range("O4").Value = 0
$ws.Range("P4").Value = "Named"
$ws.Range("Q4").Value = 63.486785924529997
$ws.Range("R4").Value = 1.76
$ws.Range("S4").Value = 0.109
$ws.Range("T4").Value = 0.0455
$ws.Range("U4").Value = 4.84
$ws.Range("V4").Value = 2.2799999999999998
$ws.Range("W4").Value = 0

# T4 needs the percent number format (unlike T2/T3 which are plain numbers)
$ws.Range("T4").NumberFormat = $ws.Range("S4").NumberFormat()
